$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("B124").Value = "Kalpataru Power"
$v = $ws.Range("B124").Value()
Write-Host "VALUE:[$v]"
